# T460 - joi 7 iulie 2022, 09:44:55 +0300
# Update the "foaie de parcurs" (trip log) for B-151-VGT, iunie 2022:
# new starting km, revised daily km/destination/purpose entries, and
# recomputed totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Km initiali (starting odometer reading)
$ws.Range("B12").Value = 263351

$ws.Range("B15").Value = 156
$ws.Range("C15").Value = "Cluj-Zalau"
$ws.Range("D15").Value = "Interes Serviciu"

$ws.Range("B16").Value = 152
$ws.Range("C16").Value = "Cluj-Cmp. Turzii"
$ws.Range("D16").Value = "Interes Serviciu"

$ws.Range("B19").Value = 421
$ws.Range("C19").Value = "Cluj-Satu-Mare"
$ws.Range("D19").Value = "Interes Serviciu"

$ws.Range("B20").Value = 356
$ws.Range("C20").Value = "Cluj-Baia-Mare"

$ws.Range("B21").Value = 121
$ws.Range("C21").Value = "Cluj-Turda"

$ws.Range("B22").Value = 156
$ws.Range("C22").Value = "Cluj-Zalau"

$ws.Range("B23").Value = 356
$ws.Range("C23").Value = "Cluj-Baia-Mare"

$ws.Range("B27").Value = 121
$ws.Range("C27").Value = "Cluj-Turda"
$ws.Range("D27").Value = "Interes Serviciu"

$ws.Range("B28").Value = 30
$ws.Range("C28").Value = "Acasa-Birou"
$ws.Range("D28").Value = " "

$ws.Range("B29").Value = 152
$ws.Range("C29").Value = "Cluj-Cmp. Turzii"

$ws.Range("B30").Value = 257
$ws.Range("C30").Value = "Cluj-Bistrita"

$ws.Range("B33").Value = 47
$ws.Range("C33").Value = "Cluj-Cluj"
$ws.Range("D33").Value = "Interes Serviciu"

$ws.Range("B34").Value = 421
$ws.Range("C34").Value = "Cluj-Satu-Mare"

$ws.Range("B36").Value = 421
$ws.Range("C36").Value = "Cluj-Satu-Mare"
$ws.Range("D36").Value = "Interes Serviciu"

$ws.Range("B37").Value = 30
$ws.Range("C37").Value = "Acasa-Birou"
$ws.Range("D37").Value = " "

$ws.Range("B40").Value = 30
$ws.Range("C40").Value = "Acasa-Birou"
$ws.Range("D40").Value = " "

$ws.Range("B41").Value = 30
$ws.Range("C41").Value = "Acasa-Birou"
$ws.Range("D41").Value = " "

$ws.Range("B42").Value = 421
$ws.Range("C42").Value = "Cluj-Satu-Mare"
$ws.Range("D42").Value = "Interes Serviciu"

$ws.Range("B43").Value = 257
$ws.Range("C43").Value = "Cluj-Bistrita"

# Totals: "Km parcursi:" (sum of daily km) and "Total" (Km initiali + Km parcursi)
$ws.Range("B44").Value = 3965

$ws.Range("B45").Value = 267316
